$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Price" (column D) values - force text storage to preserve exact formatting
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '52.020.88'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.785.44'
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '358.92'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.25'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.560'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.04'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0851'
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.43'
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.225.43'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.803.96'
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.910.55'
$ws.Range("D18").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.08'
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '272.59'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.68'
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.64'
$ws.Range("D26").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0465'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '51.31'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.04'
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.47'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0840'
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.20'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.28'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.01'
$ws.Range("D40").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '125.22'
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.84'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.068.57'
$ws.Range("D46").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.72'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.930'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.97'
$ws.Range("D51").Style = "Normal"

# Update "Volume(1h)" (column E) values
$ws.Range("E2").Value = '  -0.41%  '
$ws.Range("E3").Value = '  -1.52%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("E5").Value = '  +0.65%  '
$ws.Range("E6").Value = '  -4.32%  '
$ws.Range("E7").Value = '  +2.21%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  -2.51%  '
$ws.Range("E10").Value = '  -4.54%  '
$ws.Range("E11").Value = '  +0.18%  '
$ws.Range("E12").Value = '  +0.63%  '
$ws.Range("E13").Value = '  -3.08%  '
$ws.Range("E14").Value = '  -1.98%  '
$ws.Range("E15").Value = '  -1.69%  '
$ws.Range("E16").Value = '  -1.29%  '
$ws.Range("E17").Value = '  +4.18%  '
$ws.Range("E18").Value = '  -0.60%  '
$ws.Range("E19").Value = '  +1.06%  '
$ws.Range("E20").Value = '  -0.98%  '
$ws.Range("E21").Value = '  -4.81%  '
$ws.Range("E22").Value = '  -1.29%  '
$ws.Range("E23").Value = '  +0.53%  '
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("E25").Value = '  -1.83%  '
$ws.Range("E26").Value = '  -0.42%  '
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("E28").Value = '  -0.88%  '
$ws.Range("E29").Value = '  -1.22%  '
$ws.Range("E30").Value = '  +1.82%  '
$ws.Range("E31").Value = '  +5.45%  '
$ws.Range("E32").Value = '  +1.35%  '
$ws.Range("E33").Value = '  +0.98%  '
$ws.Range("E34").Value = '  -2.98%  '
$ws.Range("E35").Value = '  +11.61%  '
$ws.Range("E36").Value = '  +0.91%  '
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("E38").Value = '  +0.15%  '
$ws.Range("E39").Value = '  -1.35%  '
$ws.Range("E40").Value = '  -4.27%  '
$ws.Range("E41").Value = '  -0.21%  '
$ws.Range("E42").Value = '  -0.36%  '
$ws.Range("E43").Value = '  -1.77%  '
$ws.Range("E44").Value = '  -2.19%  '
$ws.Range("E45").Value = '  -7.10%  '
$ws.Range("E46").Value = '  +1.26%  '
$ws.Range("E47").Value = '  -3.30%  '
$ws.Range("E48").Value = '  -0.38%  '
$ws.Range("E49").Value = '  +0.11%  '
$ws.Range("E50").Value = '  -4.56%  '
$ws.Range("E51").Value = '  +0.92%  '
